# Update countries & provincias Spain
# Applies the 25-Abril-2020 08:52 -> 09:22 data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1) Footer timestamp string update.
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 09:22"

# 2) Lituania overtakes Bosnia y Herzegovina in the ranking (rows are sorted
#    descending by "Casos totales"/column B), so the two countries swap rows
#    76 and 77 while row 76 also carries Lituania's freshly updated figures
#    and row 77 keeps Bosnia y Herzegovina's previous (unchanged) figures.
$ws.Range("A76").Value = "Lituania"
$ws.Range("B76").Value = 1426
$ws.Range("C76").Value = 16
$ws.Range("D76").Value = 460
$ws.Range("E76").Value = 925
$ws.Range("F76").Value = 17
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 41

$ws.Range("A77").Value = "Bosnia y Herzegovina"
$ws.Range("B77").Value = 1421
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 538
$ws.Range("E77").Value = 828
$ws.Range("F77").Value = 4
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 55

# 3) Plain data refreshes for the rest of the updated countries/provinces.
# Polonia (row 33): Casos activos / Recuperados.
$ws.Range("D33").Value = 2126
$ws.Range("E33").Value = 8272

# Uzbekistan (row 67): Casos totales / Nuevos casos / Recuperados.
$ws.Range("B67").Value = 1836
$ws.Range("C67").Value = 32
$ws.Range("E67").Value = 1207

# Letonia (row 93): Casos totales / Nuevos casos / Recuperados.
$ws.Range("B93").Value = 804
$ws.Range("C93").Value = 20
$ws.Range("E93").Value = 525

# Malta (row 108): Recuperados / Casos criticos / Muertes hoy.
$ws.Range("E108").Value = 220
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 4

# Sri Lanka (row 112): Casos activos / Recuperados.
$ws.Range("D112").Value = 116
$ws.Range("E112").Value = 297
